$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param($ws, $addr, [string]$val)
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-CellText $ws "D2" "28.003.23"
Set-CellText $ws "D3" "1.890.56"
Set-CellText $ws "E3" "  -3.76%  "
Set-CellText $ws "D4" "1.002"
Set-CellText $ws "E4" "  -1.08%  "
Set-CellText $ws "D5" "326.53"
Set-CellText $ws "E5" "  +0.90%  "
Set-CellText $ws "D6" "1.003"
Set-CellText $ws "E6" "  -0.91%  "
Set-CellText $ws "D7" "0.4579"
Set-CellText $ws "E7" "  -4.19%  "
Set-CellText $ws "D8" "0.3936"
Set-CellText $ws "E8" "  -2.56%  "
Set-CellText $ws "D9" "51.48"
Set-CellText $ws "E9" "  -4.57%  "
Set-CellText $ws "D10" "0.08215"
Set-CellText $ws "E10" "  -3.13%  "
Set-CellText $ws "D11" "1.038"
Set-CellText $ws "E11" "  -2.30%  "
Set-CellText $ws "D12" "21.65"
Set-CellText $ws "E12" "  -3.52%  "
Set-CellText $ws "D13" "1.892.12"
Set-CellText $ws "E13" "  -3.00%  "
Set-CellText $ws "D14" "7.328"
Set-CellText $ws "E14" "  -4.35%  "
Set-CellText $ws "D15" "5.987"
Set-CellText $ws "E15" "  -4.07%  "
Set-CellText $ws "D16" "1.004"
Set-CellText $ws "E16" "  -0.97%  "
Set-CellText $ws "D17" "89.51"
Set-CellText $ws "E17" "  -0.42%  "
Set-CellText $ws "D18" "0.00001059"
Set-CellText $ws "E18" "  -1.12%  "
Set-CellText $ws "D19" "0.06569"
Set-CellText $ws "E19" "  -0.45%  "
Set-CellText $ws "D20" "17.59"
Set-CellText $ws "E20" "  -5.81%  "
Set-CellText $ws "E21" "  -0.79%  "
Set-CellText $ws "D22" "5.649"
Set-CellText $ws "E22" "  -2.41%  "
Set-CellText $ws "D23" "28.005.30"
Set-CellText $ws "E23" "  -2.87%  "
Set-CellText $ws "D24" "11.11"
Set-CellText $ws "E24" "  -3.82%  "
Set-CellText $ws "E25" "  +0.65%  "
Set-CellText $ws "D26" "2.146.83"
Set-CellText $ws "E26" "  -1.83%  "
Set-CellText $ws "D27" "153.98"
Set-CellText $ws "E27" "  -0.38%  "
Set-CellText $ws "D28" "19.91"
Set-CellText $ws "E28" "  -1.59%  "
Set-CellText $ws "D29" "2.110"
Set-CellText $ws "E29" "  -2.25%  "
Set-CellText $ws "D30" "5.679"
Set-CellText $ws "E30" "  -4.80%  "
Set-CellText $ws "D31" "124.41"
Set-CellText $ws "E31" "  +0.10%  "
Set-CellText $ws "D32" "0.09547"
Set-CellText $ws "E32" "  -0.66%  "
Set-CellText $ws "D33" "0.9587"
Set-CellText $ws "E33" "  -4.94%  "
Set-CellText $ws "D34" "1.480"
Set-CellText $ws "E34" "  +1.25%  "
Set-CellText $ws "E35" "  -1.32%  "
Set-CellText $ws "D36" "5.472"
Set-CellText $ws "E36" "  -3.81%  "
Set-CellText $ws "B37" "VeChain"
Set-CellText $ws "C37" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-CellText $ws "D37" "0.02282"
Set-CellText $ws "E37" "  -3.35%  "
Set-CellText $ws "B38" "TrustWalletToken"
Set-CellText $ws "C38" "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-CellText $ws "D38" "1.252"
Set-CellText $ws "E38" "  -1.21%  "
Set-CellText $ws "D39" "8.666"
Set-CellText $ws "E39" "  -1.14%  "
Set-CellText $ws "D40" "0.06110"
Set-CellText $ws "E40" "  -1.67%  "
Set-CellText $ws "D41" "0.6106"
Set-CellText $ws "E41" "  -2.12%  "
Set-CellText $ws "D42" "1.002"
Set-CellText $ws "E42" "  -0.85%  "
Set-CellText $ws "D43" "10.74"
Set-CellText $ws "E43" "  -3.51%  "
Set-CellText $ws "D44" "0.1888"
Set-CellText $ws "E44" "  -1.71%  "
Set-CellText $ws "D45" "1.309"
Set-CellText $ws "E45" "  -2.56%  "
Set-CellText $ws "D46" "0.5814"
Set-CellText $ws "E46" "  -2.57%  "
Set-CellText $ws "D47" "12.65"
Set-CellText $ws "E47" "  -2.35%  "
Set-CellText $ws "D48" "1.993"
Set-CellText $ws "E48" "  -4.05%  "
Set-CellText $ws "D49" "3.425"
Set-CellText $ws "E49" "  -0.19%  "
Set-CellText $ws "D50" "0.06881"
Set-CellText $ws "E50" "  +0.64%  "
Set-CellText $ws "D51" "110.36"
